# Add a new data block (two rows) to the bottom of the table on sheet1,
# mirroring the formatting of the previous block (rows 12:13) and
# appending the corresponding shared strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Duplicate the formatting of the last existing block (rows 12:13)
#        into the two new rows (14:15). This carries over the cell
#        styles (borders/wrap/font) used for that repeating 2-row block. ---
$ws.Range("A12:E13").Copy() | Out-Null
$ws.Range("A14:E15").PasteSpecial(-4122) | Out-Null

# --- 2. Fill in the new values (shared strings get de-duplicated/created
#        automatically by the engine when a string Value is assigned). ---

# Row 14 (first row of the new block)
$ws.Cells.Item(14,1).Value = "SCRIPT/T01P01A/um1103.ssb"
$ws.Cells.Item(14,2).Value = 201
$ws.Cells.Item(14,3).Value = " Welcome back![K] Was the\nexpedition fun?"
$ws.Cells.Item(14,4).Value = " С возвращением![K] Ну, что было\nв экспедиции?"
$ws.Cells.Item(14,5).Value = " Ò âïèâñàþåîéåí![K] Îô, œóï áúìï\nâ üëòðåäéøéé?"

# Row 15 (second row of the new block, no filename in column A)
$ws.Cells.Item(15,2).Value = 204
$ws.Cells.Item(15,3).Value = " Oh? You never discovered\nanything?[K] Oh... That\'s too bad."
$ws.Cells.Item(15,4).Value = " Ой? Вы ничего не нашли?[K]\nОх... Очень жаль."
$ws.Cells.Item(15,5).Value = " Ïê? Âú îéœåãï îå îàšìé?[K]\nÏö... Ïœåîû çàìû."

# --- 3. Row heights matching the new content's wrap (same as other
#        2-line / 1-line rows in the sheet). ---
$ws.Rows(14).RowHeight = 43.2
$ws.Rows(15).RowHeight = 21.6

# --- 4. Update the view: scroll down a bit and move the selection,
#        like the author did while adding this block. ---
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("C21").Select()
